# "10Th - MB for single stock and added new group"
#
# This MarketBeat-rank style report keeps one column per reporting date
# (most-recent date first, immediately to the right of the firm-name
# column A) and one row per covering analyst/firm. This update:
#   1. Inserts three new date columns (Jun_27, Jun_26, Jun_26) in front
#      of the existing date columns, pushing the previous Jun_17 / Jun_15
#      / Jun_13 / Jun_10 columns three places to the right.
#   2. Back-fills the new columns with the default "UN" (no action)
#      marker for every existing firm row.
#   3. Appends two new firm rows for the newly covered groups
#      ("Benchmark" and "Evercore ISI"), each pre-filled with "UN" in
#      the three newest date columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 3 new date columns before column B -----------------------
$ws.Range("B1:D1").EntireColumn.Insert()

# --- 2. New header dates for the inserted columns ------------------------
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# --- 3. Default "UN" marker for all existing firm rows (2-27) -----------
$ws.Range("B2:D27").Value = "UN"

# --- 4. Two newly covered analyst groups, appended as new rows ----------
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28:D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29:D29").Value = "UN"
